$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 386.4186
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 386.4186
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1159.2558
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1495.2558
$ws.Range("H70").Value = 7291210.5
$ws.Range("I70").Value = 16767565
$ws.Range("J70").Value = 1707.1538
$ws.Range("K70").Value = 50302695
$ws.Range("L70").Value = 5121.4614
$ws.Range("M70").Value = -50302425
$ws.Range("N70").Value = -5661.4614
$ws.Range("H73").Value = 7291210.5
$ws.Range("I73").Value = 16767565
$ws.Range("J73").Value = 1707.1538
$ws.Range("K73").Value = 50302695
$ws.Range("L73").Value = 5121.4614
$ws.Range("M73").Value = -50301759
$ws.Range("N73").Value = -6993.4614
$ws.Range("H86").Value = 8000
$ws.Range("I86").Value = 7500
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 7500
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -6377
$ws.Range("N86").Value = -10746
$ws.Range("H89").Value = 8000
$ws.Range("I89").Value = 7500
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 37500
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -31884
$ws.Range("N89").Value = -53732
$ws.Range("H125").Value = 7392.5835
$ws.Range("I125").Value = 5233
$ws.Range("J125").Value = 8472.375
$ws.Range("K125").Value = 47097
$ws.Range("L125").Value = 76251.375
$ws.Range("M125").Value = -44637
$ws.Range("N125").Value = -81171.375
$ws.Range("H132").Value = 2016.5428
$ws.Range("I132").Value = 1702.5172
$ws.Range("J132").Value = 3534.3333
$ws.Range("K132").Value = 5107.5516
$ws.Range("L132").Value = 10602.9999
$ws.Range("M132").Value = -2577.5516
$ws.Range("N132").Value = -15662.9999
$ws.Range("H135").Value = 581.8148
$ws.Range("I135").Value = 412.91666
$ws.Range("K135").Value = 3716.24994
$ws.Range("M135").Value = -1181.24994
$ws.Range("H137").Value = 2527
$ws.Range("I137").Value = 1980
$ws.Range("J137").Value = 3200.2307
$ws.Range("K137").Value = 5940
$ws.Range("L137").Value = 9600.6921
$ws.Range("M137").Value = -3390
$ws.Range("N137").Value = -14700.6921
$ws.Range("H138").Value = 3366.5
$ws.Range("I138").Value = 1303.4445
$ws.Range("J138").Value = 6019
$ws.Range("K138").Value = 3910.3335
$ws.Range("L138").Value = 18057
$ws.Range("M138").Value = 1229.6665
$ws.Range("N138").Value = -28337

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 2914.8572
$ws.Range("I5").Value = 2914.8572
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2914.8572
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2802.8572
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 23209.656
$ws.Range("I32").Value = 29616.043
$ws.Range("J32").Value = 6837.778
$ws.Range("K32").Value = 29616.043
$ws.Range("L32").Value = 6837.778
$ws.Range("M32").Value = -29329.043
$ws.Range("N32").Value = -7411.778
$ws.Range("H122").Value = 4326.727
$ws.Range("I122").Value = 4842.385
$ws.Range("J122").Value = 2411.4285
$ws.Range("K122").Value = 14527.155
$ws.Range("L122").Value = 7234.2855
$ws.Range("M122").Value = -12077.155
$ws.Range("N122").Value = -12134.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2914.8572
$ws.Range("I4").Value = 2914.8572
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2914.8572
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2799.8572
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 8437.5
$ws.Range("I5").Value = 7200
$ws.Range("J5").Value = 10500
$ws.Range("K5").Value = 7200
$ws.Range("L5").Value = 10500
$ws.Range("M5").Value = -7087
$ws.Range("N5").Value = -10726
$ws.Range("H134").Value = 2062.5264
$ws.Range("I134").Value = 2155.6924
$ws.Range("J134").Value = 1860.6666
$ws.Range("K134").Value = 6467.0772
$ws.Range("L134").Value = 5581.9998
$ws.Range("M134").Value = -3932.0772
$ws.Range("N134").Value = -10651.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 112.22222
$ws.Range("I7").Value = 47.5
$ws.Range("J7").Value = 164
$ws.Range("K7").Value = 47.5
$ws.Range("L7").Value = 164
$ws.Range("M7").Value = 65.5
$ws.Range("N7").Value = -390
$ws.Range("H22").Value = 273.85
$ws.Range("I22").Value = 248.70589
$ws.Range("K22").Value = 248.70589
$ws.Range("M22").Value = 101.29411

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 765.0678
$ws.Range("I68").Value = 565.5294
$ws.Range("J68").Value = 1036.44
$ws.Range("K68").Value = 1696.5882
$ws.Range("L68").Value = 3109.32
$ws.Range("M68").Value = -885.5882000000001
$ws.Range("N68").Value = -4731.32
$ws.Range("H71").Value = 765.0678
$ws.Range("I71").Value = 565.5294
$ws.Range("J71").Value = 1036.44
$ws.Range("K71").Value = 5089.7646
$ws.Range("L71").Value = 9327.960000000001
$ws.Range("M71").Value = -1033.7646
$ws.Range("N71").Value = -17439.96
$ws.Range("H75").Value = 11750
$ws.Range("J75").Value = 11750
$ws.Range("L75").Value = 35250
$ws.Range("N75").Value = -37246
$ws.Range("H78").Value = 11750
$ws.Range("J78").Value = 11750
$ws.Range("L78").Value = 105750
$ws.Range("N78").Value = -115734
$ws.Range("H94").Value = 3735
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3982.8572
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 11948.5716
$ws.Range("M94").Value = -5324
$ws.Range("N94").Value = -13300.5716
$ws.Range("H122").Value = 1338.5769
$ws.Range("I122").Value = 586.6667
$ws.Range("J122").Value = 1436.6522
$ws.Range("K122").Value = 5280.0003
$ws.Range("L122").Value = 12929.8698
$ws.Range("M122").Value = -2830.0003
$ws.Range("N122").Value = -17829.8698
$ws.Range("H131").Value = 2607.8484
$ws.Range("I131").Value = 605.5
$ws.Range("J131").Value = 2884.0344
$ws.Range("K131").Value = 1816.5
$ws.Range("L131").Value = 8652.1032
$ws.Range("M131").Value = 3223.5
$ws.Range("N131").Value = -18732.1032
$ws.Range("H132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 4065.8125
$ws.Range("I137").Value = 2137.1428
$ws.Range("J137").Value = 5565.8887
$ws.Range("K137").Value = 6411.428400000001
$ws.Range("L137").Value = 16697.6661
$ws.Range("M137").Value = -1311.428400000001
$ws.Range("N137").Value = -26897.6661
$ws.Range("H138").Value = 2648.923
$ws.Range("I138").Value = 906
$ws.Range("J138").Value = 3738.25
$ws.Range("K138").Value = 2718
$ws.Range("L138").Value = 11214.75
$ws.Range("M138").Value = 2422
$ws.Range("N138").Value = -21494.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.5
$ws.Range("J2").Value = 75
$ws.Range("L2").Value = 75
$ws.Range("N2").Value = -301
$ws.Range("H132").Value = 2292.9756
$ws.Range("I132").Value = 1807.1333
$ws.Range("J132").Value = 3618
$ws.Range("K132").Value = 5421.3999
$ws.Range("L132").Value = 10854
$ws.Range("M132").Value = -2891.3999
$ws.Range("N132").Value = -15914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2813.8572
$ws.Range("I61").Value = 2576.1667
$ws.Range("K61").Value = 2576.1667
$ws.Range("M61").Value = -2374.1667
$ws.Range("H81").Value = 33625
$ws.Range("J81").Value = 33625
$ws.Range("L81").Value = 33625
$ws.Range("N81").Value = -35621
$ws.Range("H82").Value = 2338.375
$ws.Range("I82").Value = 1501
$ws.Range("J82").Value = 2617.5
$ws.Range("K82").Value = 1501
$ws.Range("L82").Value = 2617.5
$ws.Range("M82").Value = -1140
$ws.Range("N82").Value = -3339.5
$ws.Range("H84").Value = 33625
$ws.Range("J84").Value = 33625
$ws.Range("L84").Value = 100875
$ws.Range("N84").Value = -110859
$ws.Range("H85").Value = 2338.375
$ws.Range("I85").Value = 1501
$ws.Range("J85").Value = 2617.5
$ws.Range("K85").Value = 1501
$ws.Range("L85").Value = 2617.5
$ws.Range("M85").Value = -253
$ws.Range("N85").Value = -5113.5
$ws.Range("H113").Value = 2813.8572
$ws.Range("I113").Value = 2576.1667
$ws.Range("K113").Value = 2576.1667
$ws.Range("M113").Value = -406.1667000000002
$ws.Range("H122").Value = 75003000
$ws.Range("I122").Value = 83335660
$ws.Range("J122").Value = 66670336
$ws.Range("K122").Value = 250006980
$ws.Range("L122").Value = 200011008
$ws.Range("M122").Value = -250004530
$ws.Range("N122").Value = -200015908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5125.375
$ws.Range("I62").Value = 5500
$ws.Range("J62").Value = 4750.75
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 4750.75
$ws.Range("M62").Value = -4876
$ws.Range("N62").Value = -5998.75
$ws.Range("H65").Value = 5125.375
$ws.Range("I65").Value = 5500
$ws.Range("J65").Value = 4750.75
$ws.Range("K65").Value = 27500
$ws.Range("L65").Value = 23753.75
$ws.Range("M65").Value = -24380
$ws.Range("N65").Value = -29993.75
$ws.Range("H68").Value = 29700
$ws.Range("J68").Value = 29700
$ws.Range("L68").Value = 29700
$ws.Range("N68").Value = -31322
$ws.Range("H71").Value = 29700
$ws.Range("J71").Value = 29700
$ws.Range("L71").Value = 89100
$ws.Range("N71").Value = -97212
